# Auto-generated Excel COM-interop script
# Applies updated market-price figures (columns H-N) across multiple leve-profit sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1898.359
$ws.Range("I17").Value = 1200
$ws.Range("J17").Value = 2208.7407
$ws.Range("K17").Value = 3600
$ws.Range("L17").Value = 6626.222099999999
$ws.Range("M17").Value = -3432
$ws.Range("N17").Value = -6962.222099999999

$ws.Range("H19").Value = 568.4706
$ws.Range("I19").Value = 649.4167
$ws.Range("J19").Value = 374.2
$ws.Range("K19").Value = 649.4167
$ws.Range("L19").Value = 374.2
$ws.Range("M19").Value = -474.4167
$ws.Range("N19").Value = -724.2

$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()

$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()

$ws.Range("H31").Value = 10000
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 10000
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 30000
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -30460

$ws.Range("H137").Value = 990222.4399999999
$ws.Range("I137").Value = 1827.1666
$ws.Range("J137").Value = 1505907
$ws.Range("K137").Value = 5481.4998
$ws.Range("L137").Value = 4517721
$ws.Range("M137").Value = -2931.4998
$ws.Range("N137").Value = -4522821

$ws.Range("H138").Value = 2944.9363
$ws.Range("I138").Value = 2193.4092
$ws.Range("J138").Value = 3606.28
$ws.Range("K138").Value = 6580.2276
$ws.Range("L138").Value = 10818.84
$ws.Range("M138").Value = -1440.2276
$ws.Range("N138").Value = -21098.84

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 870.8
$ws.Range("I25").Value = 888.5
$ws.Range("J25").Value = 800
$ws.Range("K25").Value = 888.5
$ws.Range("L25").Value = 800
$ws.Range("M25").Value = -486.5
$ws.Range("N25").Value = -1604

$ws.Range("H61").Value = 13995.091
$ws.Range("I61").Value = 13304.267
$ws.Range("J61").Value = 15475.429
$ws.Range("K61").Value = 13304.267
$ws.Range("L61").Value = 15475.429
$ws.Range("M61").Value = -13092.267
$ws.Range("N61").Value = -15899.429

$ws.Range("H74").Value = 2563.1282
$ws.Range("I74").Value = 2229.7637
$ws.Range("J74").Value = 3360.3044
$ws.Range("K74").Value = 2229.7637
$ws.Range("L74").Value = 3360.3044
$ws.Range("M74").Value = -1355.7637
$ws.Range("N74").Value = -5108.3044

$ws.Range("H77").Value = 2563.1282
$ws.Range("I77").Value = 2229.7637
$ws.Range("J77").Value = 3360.3044
$ws.Range("K77").Value = 11148.8185
$ws.Range("L77").Value = 16801.522
$ws.Range("M77").Value = -6780.818499999999
$ws.Range("N77").Value = -25537.522

$ws.Range("H122").Value = 2774.0908
$ws.Range("I122").Value = 2501.6667
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 7505.000100000001
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -5055.000100000001
$ws.Range("N122").Value = -16900

$ws.Range("H132").Value = 2652.4075
$ws.Range("I132").Value = 2281.0652
$ws.Range("J132").Value = 4787.625
$ws.Range("K132").Value = 6843.1956
$ws.Range("L132").Value = 14362.875
$ws.Range("M132").Value = -4313.1956
$ws.Range("N132").Value = -19422.875

$ws.Range("H136").Value = 13995.091
$ws.Range("I136").Value = 13304.267
$ws.Range("J136").Value = 15475.429
$ws.Range("K136").Value = 39912.801
$ws.Range("L136").Value = 46426.287
$ws.Range("M136").Value = -37362.801
$ws.Range("N136").Value = -51526.287

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1814.4762
$ws.Range("I20").Value = 1715.3572
$ws.Range("J20").Value = 2012.7142
$ws.Range("K20").Value = 1715.3572
$ws.Range("L20").Value = 2012.7142
$ws.Range("M20").Value = -1468.3572
$ws.Range("N20").Value = -2506.7142

$ws.Range("H86").Value = 1646.6666
$ws.Range("I86").Value = 1294.125
$ws.Range("J86").Value = 2351.75
$ws.Range("K86").Value = 1294.125
$ws.Range("L86").Value = 2351.75
$ws.Range("M86").Value = -171.125
$ws.Range("N86").Value = -4597.75

$ws.Range("H89").Value = 1646.6666
$ws.Range("I89").Value = 1294.125
$ws.Range("J89").Value = 2351.75
$ws.Range("K89").Value = 6470.625
$ws.Range("L89").Value = 11758.75
$ws.Range("M89").Value = -854.625
$ws.Range("N89").Value = -22990.75

$ws.Range("H107").Value = 1388.1818
$ws.Range("I107").Value = 1333.0488
$ws.Range("J107").Value = 2141.6667
$ws.Range("K107").Value = 1333.0488
$ws.Range("L107").Value = 2141.6667
$ws.Range("M107").Value = 586.9512
$ws.Range("N107").Value = -5981.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4995.241
$ws.Range("I31").Value = 1857.9032
$ws.Range("J31").Value = 8597.370000000001
$ws.Range("K31").Value = 1857.9032
$ws.Range("L31").Value = 8597.370000000001
$ws.Range("M31").Value = -1562.9032
$ws.Range("N31").Value = -9187.370000000001

$ws.Range("H34").Value = 4995.241
$ws.Range("I34").Value = 1857.9032
$ws.Range("J34").Value = 8597.370000000001
$ws.Range("K34").Value = 1857.9032
$ws.Range("L34").Value = 8597.370000000001
$ws.Range("M34").Value = -1655.9032
$ws.Range("N34").Value = -9001.370000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 7083.9
$ws.Range("I3").Value = 1415.8572
$ws.Range("J3").Value = 20309.334
$ws.Range("K3").Value = 4247.571599999999
$ws.Range("L3").Value = 60928.00199999999
$ws.Range("M3").Value = -4135.571599999999
$ws.Range("N3").Value = -61152.00199999999

$ws.Range("H6").Value = 20.666666
$ws.Range("I6").Value = 20.666666
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 61.999998
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 51.000002
$ws.Range("N6").ClearContents()

$ws.Range("H12").Value = 21.153847
$ws.Range("I12").Value = 10.2
$ws.Range("J12").Value = 28
$ws.Range("K12").Value = 30.6
$ws.Range("L12").Value = 84
$ws.Range("M12").Value = 142.4
$ws.Range("N12").Value = -430

$ws.Range("H93").Value = 5210.1
$ws.Range("I93").Value = 1999
$ws.Range("J93").Value = 5566.8887
$ws.Range("K93").Value = 5997
$ws.Range("L93").Value = 16700.6661
$ws.Range("M93").Value = -4125
$ws.Range("N93").Value = -20444.6661

$ws.Range("H136").Value = 2810.9092
$ws.Range("I136").Value = 2612
$ws.Range("J136").Value = 4800
$ws.Range("K136").Value = 7836
$ws.Range("L136").Value = 14400
$ws.Range("M136").Value = -2736
$ws.Range("N136").Value = -24600

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 288844.7
$ws.Range("I132").Value = 464638.62
$ws.Range("J132").Value = 3179.5625
$ws.Range("K132").Value = 1393915.86
$ws.Range("L132").Value = 9538.6875
$ws.Range("M132").Value = -1391385.86
$ws.Range("N132").Value = -14598.6875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 16999.857
$ws.Range("I16").Value = 25999.75
$ws.Range("J16").Value = 5000
$ws.Range("K16").Value = 25999.75
$ws.Range("L16").Value = 5000
$ws.Range("M16").Value = -25829.75
$ws.Range("N16").Value = -5340

$ws.Range("H22").Value = 43817.875
$ws.Range("I22").Value = 200769.6
$ws.Range("J22").Value = 2514.7896
$ws.Range("K22").Value = 200769.6
$ws.Range("L22").Value = 2514.7896
$ws.Range("M22").Value = -200474.6
$ws.Range("N22").Value = -3104.7896

$ws.Range("H27").Value = 43817.875
$ws.Range("I27").Value = 200769.6
$ws.Range("J27").Value = 2514.7896
$ws.Range("K27").Value = 200769.6
$ws.Range("L27").Value = 2514.7896
$ws.Range("M27").Value = -200662.6
$ws.Range("N27").Value = -2728.7896

$ws.Range("H46").Value = 3247.2903
$ws.Range("I46").Value = 1983.3334
$ws.Range("J46").Value = 3382.7144
$ws.Range("K46").Value = 1983.3334
$ws.Range("L46").Value = 3382.7144
$ws.Range("M46").Value = -1795.3334
$ws.Range("N46").Value = -3758.7144

$ws.Range("H55").Value = 210.72223
$ws.Range("I55").Value = 165.91667
$ws.Range("J55").Value = 300.33334
$ws.Range("K55").Value = 165.91667
$ws.Range("L55").Value = 300.33334
$ws.Range("M55").Value = 7.083329999999989
$ws.Range("N55").Value = -646.33334

$ws.Range("H109").Value = 39000
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 39000
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 39000
$ws.Range("N109").Value = -41774

$ws.Range("H111").Value = 0
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()

$ws.Range("H132").Value = 788490.4399999999
$ws.Range("I132").Value = 1049960
$ws.Range("J132").Value = 4081.7273
$ws.Range("K132").Value = 3149880
$ws.Range("L132").Value = 12245.1819
$ws.Range("M132").Value = -3147350
$ws.Range("N132").Value = -17305.1819

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 701
$ws.Range("I81").Value = 699.4
$ws.Range("J81").Value = 703.6667
$ws.Range("K81").Value = 1398.8
$ws.Range("L81").Value = 1407.3334
$ws.Range("M81").Value = -337.8
$ws.Range("N81").Value = -3529.3334

$ws.Range("H84").Value = 701
$ws.Range("I84").Value = 699.4
$ws.Range("J84").Value = 703.6667
$ws.Range("K84").Value = 6994
$ws.Range("L84").Value = 7036.666999999999
$ws.Range("M84").Value = -1690
$ws.Range("N84").Value = -17644.667

$ws.Range("H113").Value = 809.76
$ws.Range("I113").Value = 638.0714
$ws.Range("J113").Value = 1028.2727
$ws.Range("K113").Value = 1914.2142
$ws.Range("L113").Value = 3084.8181
$ws.Range("M113").Value = 255.7857999999999
$ws.Range("N113").Value = -7424.8181

$ws.Range("H122").Value = 2656.625
$ws.Range("I122").Value = 1409.1111
$ws.Range("J122").Value = 4260.5713
$ws.Range("K122").Value = 4227.3333
$ws.Range("L122").Value = 12781.7139
$ws.Range("M122").Value = -1777.3333
$ws.Range("N122").Value = -17681.7139

$ws.Range("H126").Value = 2654.9443
$ws.Range("I126").Value = 2439.8
$ws.Range("J126").Value = 2923.875
$ws.Range("K126").Value = 7319.400000000001
$ws.Range("L126").Value = 8771.625
$ws.Range("M126").Value = -4849.400000000001
$ws.Range("N126").Value = -13711.625

$ws.Range("H132").Value = 475899.34
$ws.Range("I132").Value = 737575.0600000001
$ws.Range("J132").Value = 2870.1155
$ws.Range("K132").Value = 2212725.18
$ws.Range("L132").Value = 8610.3465
$ws.Range("M132").Value = -2210195.18
$ws.Range("N132").Value = -13670.3465

$ws.Range("H136").Value = 7603352.5
$ws.Range("I136").Value = 9271785
$ws.Range("J136").Value = 2716.4443
$ws.Range("K136").Value = 27815355
$ws.Range("L136").Value = 8149.3329
$ws.Range("M136").Value = -27812805
$ws.Range("N136").Value = -13249.3329
